# Auto-generated edit script: apply cached-value updates to Cactuar_Profits workbook
# Values come from a scheduled market-data refresh; cells are plain numeric
# literals (no formulas) in columns H-N of each sheet.
$wb = $excel.ActiveWorkbook

# --- Worksheet #1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(33, 8).Value = 417.85715  # H33: was 433.7857
$ws.Cells.Item(33, 9).Value = 449.81818  # I33: was 470.0909
$ws.Cells.Item(33, 11).Value = 449.81818  # K33: was 470.0909
$ws.Cells.Item(33, 13).Value = -220.81818  # M33: was -241.0909
$ws.Cells.Item(52, 8).Value = 13890388  # H52: was 8335681.5
$ws.Cells.Item(52, 9).Value = 13890388  # I52: was 10418391
$ws.Cells.Item(52, 10).Value = 0  # J52: was 4845
$ws.Cells.Item(52, 11).Value = 41671164  # K52: was 31255173
$ws.Cells.Item(52, 12).Value = 0  # L52: was 14535
$ws.Cells.Item(52, 13).Value = -41671004  # M52: was -31255013
$ws.Cells.Item(52, 14).ClearContents()  # N52: was -14855
$ws.Cells.Item(111, 8).Value = 5500.5386  # H111: was 5940.1665
$ws.Cells.Item(111, 9).Value = 3955.182  # I111: was 4328.2
$ws.Cells.Item(111, 11).Value = 11865.546  # K111: was 12984.6
$ws.Cells.Item(111, 13).Value = -8798.545999999998  # M111: was -9917.599999999999
$ws.Cells.Item(121, 8).Value = 4797.609  # H121: was 4834.8184
$ws.Cells.Item(121, 10).Value = 4797.609  # J121: was 4834.8184
$ws.Cells.Item(121, 12).Value = 14392.827  # L121: was 14504.4552
$ws.Cells.Item(121, 14).Value = -17886.827  # N121: was -17998.4552
$ws.Cells.Item(135, 8).Value = 4678.44  # H135: was 4682.48
$ws.Cells.Item(135, 9).Value = 1161.9375  # I135: was 1168.25
$ws.Cells.Item(135, 11).Value = 10457.4375  # K135: was 10514.25
$ws.Cells.Item(135, 13).Value = -7922.4375  # M135: was -7979.25
$ws.Cells.Item(137, 8).Value = 11114355  # H137: was 10755843
$ws.Cells.Item(137, 9).Value = 1771.7  # I137: was 1709.9524
$ws.Cells.Item(137, 11).Value = 5315.1  # K137: was 5129.857199999999
$ws.Cells.Item(137, 13).Value = -2765.1  # M137: was -2579.857199999999

# --- Worksheet #2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 16350.702  # H32: was 17085.125
$ws.Cells.Item(32, 9).Value = 16059.84  # I32: was 16724.396
$ws.Cells.Item(32, 10).Value = 18428.285  # J32: was 19249.5
$ws.Cells.Item(32, 11).Value = 16059.84  # K32: was 16724.396
$ws.Cells.Item(32, 12).Value = 18428.285  # L32: was 19249.5
$ws.Cells.Item(32, 13).Value = -15772.84  # M32: was -16437.396
$ws.Cells.Item(32, 14).Value = -19002.285  # N32: was -19823.5
$ws.Cells.Item(61, 8).Value = 4331.02  # H61: was 4415.3267
$ws.Cells.Item(61, 9).Value = 3370.7568  # I61: was 3458.8333
$ws.Cells.Item(61, 11).Value = 3370.7568  # K61: was 3458.8333
$ws.Cells.Item(61, 13).Value = -3158.7568  # M61: was -3246.8333
$ws.Cells.Item(74, 8).Value = 5435485.5  # H74: was 5435488
$ws.Cells.Item(74, 9).Value = 6250595.5  # I74: was 6250598.5
$ws.Cells.Item(74, 11).Value = 6250595.5  # K74: was 6250598.5
$ws.Cells.Item(74, 13).Value = -6249721.5  # M74: was -6249724.5
$ws.Cells.Item(77, 8).Value = 5435485.5  # H77: was 5435488
$ws.Cells.Item(77, 9).Value = 6250595.5  # I77: was 6250598.5
$ws.Cells.Item(77, 11).Value = 31252977.5  # K77: was 31252992.5
$ws.Cells.Item(77, 13).Value = -31248609.5  # M77: was -31248624.5
$ws.Cells.Item(136, 8).Value = 4331.02  # H136: was 4415.3267
$ws.Cells.Item(136, 9).Value = 3370.7568  # I136: was 3458.8333
$ws.Cells.Item(136, 11).Value = 10112.2704  # K136: was 10376.4999
$ws.Cells.Item(136, 13).Value = -7562.270400000001  # M136: was -7826.499899999999

# --- Worksheet #3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 1655.4857  # H134: was 1752.0294
$ws.Cells.Item(134, 9).Value = 1333.3793  # I134: was 1417.0358
$ws.Cells.Item(134, 10).Value = 3212.3333  # J134: was 3315.3333
$ws.Cells.Item(134, 11).Value = 4000.1379  # K134: was 4251.107400000001
$ws.Cells.Item(134, 12).Value = 9636.999899999999  # L134: was 9945.999899999999
$ws.Cells.Item(134, 13).Value = -1465.1379  # M134: was -1716.107400000001
$ws.Cells.Item(134, 14).Value = -14706.9999  # N134: was -15015.9999

# --- Worksheet #4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(22, 8).Value = 1829.8572  # H22: was 1707.375
$ws.Cells.Item(22, 9).Value = 1607.6666  # I22: was 1418.25
$ws.Cells.Item(22, 11).Value = 1607.6666  # K22: was 1418.25
$ws.Cells.Item(22, 13).Value = -1257.6666  # M22: was -1068.25
$ws.Cells.Item(31, 8).Value = 1901.1195  # H31: was 1881.9043
$ws.Cells.Item(31, 9).Value = 1575.3928  # I31: was 1561.9651
$ws.Cells.Item(31, 11).Value = 1575.3928  # K31: was 1561.9651
$ws.Cells.Item(31, 13).Value = -1280.3928  # M31: was -1266.9651
$ws.Cells.Item(34, 8).Value = 1901.1195  # H34: was 1881.9043
$ws.Cells.Item(34, 9).Value = 1575.3928  # I34: was 1561.9651
$ws.Cells.Item(34, 11).Value = 1575.3928  # K34: was 1561.9651
$ws.Cells.Item(34, 13).Value = -1373.3928  # M34: was -1359.9651
$ws.Cells.Item(99, 8).Value = 6761.5  # H99: was 7185.1816
$ws.Cells.Item(99, 9).Value = 4705.4287  # I99: was 4833.857
$ws.Cells.Item(99, 10).Value = 9640  # J99: was 11300
$ws.Cells.Item(99, 11).Value = 4705.4287  # K99: was 4833.857
$ws.Cells.Item(99, 12).Value = 9640  # L99: was 11300
$ws.Cells.Item(99, 13).Value = -3207.4287  # M99: was -3335.857
$ws.Cells.Item(99, 14).Value = -12636  # N99: was -14296
$ws.Cells.Item(126, 8).Value = 6761.5  # H126: was 7185.1816
$ws.Cells.Item(126, 9).Value = 4705.4287  # I126: was 4833.857
$ws.Cells.Item(126, 10).Value = 9640  # J126: was 11300
$ws.Cells.Item(126, 11).Value = 14116.2861  # K126: was 14501.571
$ws.Cells.Item(126, 12).Value = 28920  # L126: was 33900
$ws.Cells.Item(126, 13).Value = -11646.2861  # M126: was -12031.571
$ws.Cells.Item(126, 14).Value = -33860  # N126: was -38840
$ws.Cells.Item(132, 8).Value = 12826491  # H132: was 9263751
$ws.Cells.Item(132, 9).Value = 19609704  # I132: was 12347082
$ws.Cells.Item(132, 11).Value = 58829112  # K132: was 37041246
$ws.Cells.Item(132, 13).Value = -58826582  # M132: was -37038716
$ws.Cells.Item(134, 8).Value = 2943.5938  # H134: was 3006.258
$ws.Cells.Item(134, 9).Value = 2024.3529  # I134: was 2088.3125
$ws.Cells.Item(134, 11).Value = 6073.0587  # K134: was 6264.9375
$ws.Cells.Item(134, 13).Value = -3538.0587  # M134: was -3729.9375

# --- Worksheet #5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(51, 8).Value = 399.5  # H51: was 421
$ws.Cells.Item(51, 9).Value = 399.5  # I51: was 421
$ws.Cells.Item(51, 11).Value = 1198.5  # K51: was 1263
$ws.Cells.Item(51, 13).Value = -738.5  # M51: was -803
$ws.Cells.Item(113, 8).Value = 650.7778  # H113: was 663
$ws.Cells.Item(113, 10).Value = 607.25  # J113: was 615
$ws.Cells.Item(113, 12).Value = 1821.75  # L113: was 1845
$ws.Cells.Item(113, 14).Value = -6161.75  # N113: was -6185
$ws.Cells.Item(117, 8).Value = 4005.6365  # H117: was 4213
$ws.Cells.Item(117, 10).Value = 4295.25  # J117: was 4632.857
$ws.Cells.Item(117, 12).Value = 12885.75  # L117: was 13898.571
$ws.Cells.Item(117, 14).Value = -19769.75  # N117: was -20782.571
$ws.Cells.Item(127, 8).Value = 0  # H127: was 1000
$ws.Cells.Item(127, 10).Value = 0  # J127: was 1000
$ws.Cells.Item(127, 12).Value = 0  # L127: was 3000
$ws.Cells.Item(127, 14).ClearContents()  # N127: was -12920
$ws.Cells.Item(131, 8).Value = 22007544  # H131: was 18863738
$ws.Cells.Item(131, 9).Value = 33334286  # I131: was 18519380
$ws.Cells.Item(131, 10).Value = 20476902  # J131: was 18941218
$ws.Cells.Item(131, 11).Value = 100002858  # K131: was 55558140
$ws.Cells.Item(131, 12).Value = 61430706  # L131: was 56823654
$ws.Cells.Item(131, 13).Value = -99997818  # M131: was -55553100
$ws.Cells.Item(131, 14).Value = -61440786  # N131: was -56833734

# --- Worksheet #6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(62, 8).Value = 0  # H62: was 59000
$ws.Cells.Item(62, 10).Value = 0  # J62: was 59000
$ws.Cells.Item(62, 12).Value = 0  # L62: was 59000
$ws.Cells.Item(62, 14).ClearContents()  # N62: was -60372
$ws.Cells.Item(65, 8).Value = 0  # H65: was 59000
$ws.Cells.Item(65, 10).Value = 0  # J65: was 59000
$ws.Cells.Item(65, 12).Value = 0  # L65: was 177000
$ws.Cells.Item(65, 14).ClearContents()  # N65: was -183864
$ws.Cells.Item(70, 8).Value = 7673.1665  # H70: was 7799.5
$ws.Cells.Item(70, 9).Value = 7263.684  # I70: was 7463.1577
$ws.Cells.Item(70, 11).Value = 7263.684  # K70: was 7463.1577
$ws.Cells.Item(70, 13).Value = -6993.684  # M70: was -7193.1577
$ws.Cells.Item(73, 8).Value = 7673.1665  # H73: was 7799.5
$ws.Cells.Item(73, 9).Value = 7263.684  # I73: was 7463.1577
$ws.Cells.Item(73, 11).Value = 7263.684  # K73: was 7463.1577
$ws.Cells.Item(73, 13).Value = -6327.684  # M73: was -6527.1577
$ws.Cells.Item(132, 8).Value = 4916.1665  # H132: was 4499.4165
$ws.Cells.Item(132, 9).Value = 4299.4  # I132: was 4363
$ws.Cells.Item(132, 10).Value = 8000  # J132: was 6000
$ws.Cells.Item(132, 11).Value = 12898.2  # K132: was 13089
$ws.Cells.Item(132, 12).Value = 24000  # L132: was 18000
$ws.Cells.Item(132, 13).Value = -10368.2  # M132: was -10559
$ws.Cells.Item(132, 14).Value = -29060  # N132: was -23060

# --- Worksheet #7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 4074.1482  # H7: was 4050.9268
$ws.Cells.Item(7, 9).Value = 3449.8704  # I7: was 3440.3035
$ws.Cells.Item(7, 10).Value = 5322.7036  # J7: was 5366.115
$ws.Cells.Item(7, 11).Value = 3449.8704  # K7: was 3440.3035
$ws.Cells.Item(7, 12).Value = 5322.7036  # L7: was 5366.115
$ws.Cells.Item(7, 13).Value = -3337.8704  # M7: was -3328.3035
$ws.Cells.Item(7, 14).Value = -5546.7036  # N7: was -5590.115
$ws.Cells.Item(16, 8).Value = 2072.9524  # H16: was 1981
$ws.Cells.Item(16, 10).Value = 4162.1665  # J16: was 3574.7144
$ws.Cells.Item(16, 12).Value = 4162.1665  # L16: was 3574.7144
$ws.Cells.Item(16, 14).Value = -4502.1665  # N16: was -3914.7144
$ws.Cells.Item(22, 8).Value = 734.9167  # H22: was 735.1739
$ws.Cells.Item(22, 9).Value = 706.8333  # I22: was 657.9231
$ws.Cells.Item(22, 10).Value = 763  # J22: was 835.6
$ws.Cells.Item(22, 11).Value = 706.8333  # K22: was 657.9231
$ws.Cells.Item(22, 12).Value = 763  # L22: was 835.6
$ws.Cells.Item(22, 13).Value = -411.8333  # M22: was -362.9231
$ws.Cells.Item(22, 14).Value = -1353  # N22: was -1425.6
$ws.Cells.Item(27, 8).Value = 734.9167  # H27: was 735.1739
$ws.Cells.Item(27, 9).Value = 706.8333  # I27: was 657.9231
$ws.Cells.Item(27, 10).Value = 763  # J27: was 835.6
$ws.Cells.Item(27, 11).Value = 706.8333  # K27: was 657.9231
$ws.Cells.Item(27, 12).Value = 763  # L27: was 835.6
$ws.Cells.Item(27, 13).Value = -599.8333  # M27: was -550.9231
$ws.Cells.Item(27, 14).Value = -977  # N27: was -1049.6
$ws.Cells.Item(126, 8).Value = 4074.1482  # H126: was 4050.9268
$ws.Cells.Item(126, 9).Value = 3449.8704  # I126: was 3440.3035
$ws.Cells.Item(126, 10).Value = 5322.7036  # J126: was 5366.115
$ws.Cells.Item(126, 11).Value = 10349.6112  # K126: was 10320.9105
$ws.Cells.Item(126, 12).Value = 15968.1108  # L126: was 16098.345
$ws.Cells.Item(126, 13).Value = -7879.611199999999  # M126: was -7850.9105
$ws.Cells.Item(126, 14).Value = -20908.1108  # N126: was -21038.345
$ws.Cells.Item(132, 8).Value = 2698.59  # H132: was 2680.889
$ws.Cells.Item(132, 9).Value = 2662.4102  # I132: was 2639.182
$ws.Cells.Item(132, 11).Value = 7987.230599999999  # K132: was 7917.545999999999
$ws.Cells.Item(132, 13).Value = -5457.230599999999  # M132: was -5387.545999999999

# --- Worksheet #8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(122, 8).Value = 2287  # H122: was 2364.7856
$ws.Cells.Item(122, 9).Value = 1345.421  # I122: was 1414.1111
$ws.Cells.Item(122, 11).Value = 4036.263  # K122: was 4242.3333
$ws.Cells.Item(122, 13).Value = -1586.263  # M122: was -1792.3333
